# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (Exhibitions) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3582
$ws1.Range("F5").Value = 2214
$ws1.Range("F6").Value = 432
$ws1.Range("F7").Value = 1
$ws1.Range("F9").Value = 77
$ws1.Range("F10").Value = 66
$ws1.Range("F11").Value = 1322
$ws1.Range("F13").Value = 1884

# --- Sheet "全部类型" ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F4").Value = 3582
$ws2.Range("F5").Value = 2214
$ws2.Range("F6").Value = 432
$ws2.Range("F7").Value = 1
$ws2.Range("F10").Value = 77
$ws2.Range("F11").Value = 66
$ws2.Range("F14").Value = 1322
$ws2.Range("F16").Value = 1884
